# Update the lattice-multiplication exercise table: each of the 15
# cells keeps its existing 5-line layout (problem / top factors /
# "----" rule / two leading-digit rows), only the numbers change.
# Word represents the <w:br/> line breaks inside a cell's Range.Text
# as vertical-tab (chr 11) characters, so we rebuild each cell's text
# using that separator - this regenerates the same run/break structure
# on save.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

# Row-major list of the 15 new cell contents (5 lines each).
$newCells = @(
    @("12 x 97", "  9    7", "  ----", "1|    |", "2|    |"),
    @("49 x 37", "  3    7", "  ----", "4|    |", "9|    |"),
    @("87 x 83", "  8    3", "  ----", "8|    |", "7|    |"),

    @("87 x 91", "  9    1", "  ----", "8|    |", "7|    |"),
    @("75 x 44", "  4    4", "  ----", "7|    |", "5|    |"),
    @("89 x 17", "  1    7", "  ----", "8|    |", "9|    |"),

    @("81 x 38", "  3    8", "  ----", "8|    |", "1|    |"),
    @("56 x 54", "  5    4", "  ----", "5|    |", "6|    |"),
    @("19 x 17", "  1    7", "  ----", "1|    |", "9|    |"),

    @("92 x 37", "  3    7", "  ----", "9|    |", "2|    |"),
    @("16 x 86", "  8    6", "  ----", "1|    |", "6|    |"),
    @("65 x 67", "  6    7", "  ----", "6|    |", "5|    |"),

    @("68 x 34", "  3    4", "  ----", "6|    |", "8|    |"),
    @("99 x 37", "  3    7", "  ----", "9|    |", "9|    |"),
    @("87 x 44", "  4    4", "  ----", "8|    |", "7|    |")
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $lines = $newCells[$idx]
        $newText = [string]::Join($nl, $lines)
        $t.Cell($r, $c).Range.Text = $newText
        $idx++
    }
}
